# Updated cryptos list on Tue Dec 19 02:42:06 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) figures on the cryptos sheet and
# reorders the InjectiveProtocol / Monero rows (31-32) to reflect the
# updated rankings, together with their new Price / Volume(1h) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.178.71"
$ws.Range("E2").Value = "  +5.47%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.236.69"
$ws.Range("E3").Value = "  +3.12%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.34%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.70%  "

# Row 7 - Solana
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.29%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.25%  "

# Row 9 - Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.615"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.70%  "

# Row 10 - Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.65%  "

# Row 11 - Dogecoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.34%  "

# Row 12 - OKB
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.97%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.99"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.01%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +0.62%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.574.09"
$ws.Range("E15").Value = "  +3.01%  "

# Row 16 - Chainlink
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.93%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.235.79"
$ws.Range("E17").Value = "  +3.47%  "

# Row 18 - Polygon
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.817"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.12%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "43.049.35"
$ws.Range("E19").Value = "  +5.27%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  +3.71%  "

# Row 21 - Litecoin
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.26%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.49%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.51%  "

# Row 24 - BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "230.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.81%  "

# Row 25 - ImmutableX
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.36%  "

# Row 26 - Dai
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27 - Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "

# Row 28 - WEMIXToken
$ws.Range("E28").Value = "  -5.63%  "

# Row 29 - PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.21%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +0.32%  "

# Rows 31-32 - InjectiveProtocol overtakes Monero in the rankings, so the
# two rows swap places (name/link) and pick up refreshed Price / Volume.
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +23.32%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "174.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.18%  "

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.11%  "

# Row 34 - Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0798"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.40%  "

# Row 35 - Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.65%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  +1.74%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  +7.42%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +5.60%  "

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0331"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +15.33%  "

# Row 40 - Celestia
$ws.Range("E40").Value = "  +7.28%  "

# Row 41 - LidoDAOToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.81%  "

# Row 42 - THORChain
$ws.Range("E42").Value = "  +3.29%  "

# Row 43 - Algorand
$ws.Range("E43").Value = "  +5.56%  "

# Row 44 - MultiversX
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "

# Row 45 - Aave
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.35%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  +3.35%  "

# Row 47 - Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.53%  "

# Row 48 - WOONetwork
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.458"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +23.95%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +2.83%  "

# Row 50 - NEARProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.55%  "

# Row 51 - TrustWalletToken
$ws.Range("E51").Value = "  +2.69%  "
